$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 513; existing rows 513:613 shift down to 514:614.
$ws.Rows(513).Insert()

# Populate the newly inserted row 513 with the new "Ciruela / Lemon / Primera"
# price record (same market/date/variety context as its neighbours, new
# quality/volume/price/unit/ratio figures).
$ws.Range("A513").Value = 6
$ws.Range("B513").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C513").Value = "Metropolitana"
$ws.Range("D513").Value = 44222
$ws.Range("E513").Value = 13
$ws.Range("F513").Value = "Fruta"
$ws.Range("G513").Value = 100103
$ws.Range("H513").Value = "Frutos de hueso (carozo)"
$ws.Range("I513").Value = 100103002
$ws.Range("J513").Value = "Ciruela"
$ws.Range("K513").Value = "Lemon"
$ws.Range("L513").Value = "Primera"
$ws.Range("M513").Value = 150
$ws.Range("N513").Value = 90000
$ws.Range("O513").Value = 90000
$ws.Range("P513").Value = 90000
$ws.Range("Q513").Value = "$/caja 15 kilos granel"
$ws.Range("R513").Value = "Región de O'Higgins"
$ws.Range("S513").Value = 6000
$ws.Range("T513").Value = 15
